# Update Volume(1h) percentage values in column E (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = "  -3.23%  "
    3 = "  -5.09%  "
    4 = "  +0.25%  "
    5 = "  -5.21%  "
    6 = "  -7.05%  "
    7 = "  +0.12%  "
    8 = "  -5.21%  "
    9 = "  -6.40%  "
    10 = "  -7.69%  "
    11 = "  -10.93%  "
    12 = "  -10.04%  "
    13 = "  -4.16%  "
    14 = "  -1.25%  "
    15 = "  -7.87%  "
    16 = "  -4.00%  "
    17 = "  -7.55%  "
    18 = "  -10.09%  "
    19 = "  -7.33%  "
    20 = "  -6.83%  "
    21 = "  -9.29%  "
    22 = "  -10.01%  "
    23 = "  -0.28%  "
    24 = "  -7.25%  "
    25 = "  -6.58%  "
    26 = "  -3.45%  "
    27 = "  +0.14%  "
    28 = "  -11.51%  "
    29 = "  +0.07%  "
    30 = "  -7.56%  "
    31 = "  -2.37%  "
    32 = "  -7.87%  "
    33 = "  -7.33%  "
    34 = "  -6.40%  "
    35 = "  -3.95%  "
    36 = "  -7.86%  "
    37 = "  -8.58%  "
    38 = "  -10.51%  "
    39 = "  -4.44%  "
    40 = "  -7.45%  "
    41 = "  -10.90%  "
    42 = "  -11.25%  "
    43 = "  +0.17%  "
    44 = "  -8.36%  "
    45 = "  -6.12%  "
    46 = "  -8.31%  "
    47 = "  -3.26%  "
    48 = "  -9.94%  "
    49 = "  -5.90%  "
    50 = "  -6.55%  "
    51 = "  -7.36%  "
}

foreach ($row in $newValues.Keys) {
    $ws.Range("E" + $row).Value = $newValues[$row]
}
